$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 13 (the lone "Docentes responsaveis" value row with no label),
# which shifts rows 14-24 up to become rows 13-23 and realigns row heights.
$ws.Range("A13:C13").EntireRow.Delete()

# Fix up cell contents that differ between the shifted-up old content and the target content.
$ws.Range("B10").Value = '5840897 - Clodoaldo Saron'
$ws.Range("C10").Value = '5840897 - Clodoaldo Saron'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("B15").Value = '01/01/2012'
$ws.Range("C15").Value = '01/01/2012'
$ws.Range("B18").Value = '5840897 - Clodoaldo Saron'
$ws.Range("C18").Value = '5840897 - Clodoaldo Saron'
$ws.Range("B19").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Range("C19").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Range("B20").Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Range("C20").Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
